$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.858.56'
$ws.Range("E2").Value = '  -1.07%  '

$ws.Range("D3").Value = '2.450.72'
$ws.Range("E3").Value = '  -2.93%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.08'
$ws.Range("E5").Value = '  +0.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.14'
$ws.Range("E6").Value = '  -1.38%  '

$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("D9").Value = '2.455.89'
$ws.Range("E9").Value = '  -2.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0976'
$ws.Range("E10").Value = '  -0.20%  '

$ws.Range("E11").Value = '  -1.87%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.97'
$ws.Range("E12").Value = '  -3.70%  '

$ws.Range("E13").Value = '  -2.44%  '

$ws.Range("D14").Value = '2.885.79'
$ws.Range("E14").Value = '  -2.86%  '

$ws.Range("D15").Value = '57.802.12'
$ws.Range("E15").Value = '  -1.10%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.69'
$ws.Range("E16").Value = '  -1.85%  '

$ws.Range("E17").Value = '  -1.53%  '

$ws.Range("D18").Value = '2.453.37'
$ws.Range("E18").Value = '  -2.85%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.29'
$ws.Range("E19").Value = '  -3.50%  '

$ws.Range("E20").Value = '  -1.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '311.26'
$ws.Range("E21").Value = '  -3.35%  '

$ws.Range("E22").Value = '  -1.20%  '

$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.89'
$ws.Range("E24").Value = '  +0.34%  '

$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '2.587.43'
$ws.Range("E25").Value = '  -1.68%  '

$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.402'
$ws.Range("E26").Value = '  -1.18%  '

$ws.Range("E27").Value = '  -0.13%  '

$ws.Range("E28").Value = '  -1.62%  '

$ws.Range("E29").Value = '  -2.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '172.97'
$ws.Range("E30").Value = '  +2.70%  '

$ws.Range("E31").Value = '  -1.99%  '

$ws.Range("E32").Value = '  -1.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.21'
$ws.Range("E33").Value = '  -1.32%  '

$ws.Range("E34").Value = '  -4.36%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.78'
$ws.Range("E37").Value = '  -2.13%  '

$ws.Range("E38").Value = '  -4.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.78'
$ws.Range("E39").Value = '  -3.36%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.815'
$ws.Range("E40").Value = '  +5.69%  '

$ws.Range("E41").Value = '  -0.72%  '

$ws.Range("E42").Value = '  -2.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.39'
$ws.Range("E43").Value = '  -1.57%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '261.86'
$ws.Range("E44").Value = '  -5.35%  '

$ws.Range("E45").Value = '  -2.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.80'
$ws.Range("E46").Value = '  -3.51%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0919'
$ws.Range("E47").Value = '  +0.28%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.43'
$ws.Range("E48").Value = '  -5.73%  '

$ws.Range("E50").Value = '  -1.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.00'
$ws.Range("E51").Value = '  -3.88%  '
